$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input")

# Fix the OrderId for row 2 (demo script value correction)
$ws.Range("Z2").NumberFormatLocal = "@"
$ws.Range("Z2").Value = "51490999"

# Fix the orderReferenceID for row 3 (was a leftover "Clone_" demo value)
$ws.Range("C3").Value = "DIR_C02"
